# repull data, push all data, mean calculation
# Update the dSF column (F) values for the rows whose computed "final minus
# initial" metric changed after the repull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 1
    4  = -1
    5  = -1
    6  = -3
    7  = 3
    8  = -1
    9  = -1
    10 = -5
    12 = -4
    13 = 2
    14 = -3
    15 = -2
    18 = -4
    19 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
